$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws4 = $wb.Worksheets.Item(4)  # 全部类型


# ---- Sheet 1 ----
$ws1.Range("F2").Value = 327
$ws1.Range("F4").Value = 611
$ws1.Range("F6").Value = 453
$ws1.Range("F7").Value = 36
$ws1.Range("F8").Value = 2102
$ws1.Range("F9").Value = 845
$ws1.Range("F10").Value = 809
$ws1.Range("F11").Value = 399
$ws1.Range("F12").Value = 63
$ws1.Range("F14").Value = 317
$ws1.Range("F16").Value = 680
$ws1.Range("F18").Value = 25
$ws1.Range("F19").Value = 1626
$ws1.Range("F20").Value = 41
$ws1.Range("F21").Value = 25
$ws1.Range("F22").Value = 23
$ws1.Range("F25").Value = 1450
$ws1.Range("F27").Value = 515
$ws1.Range("F29").Value = 566
$ws1.Range("F30").Value = 406
$ws1.Range("F31").Value = 2269
$ws1.Range("F32").Value = 379
$ws1.Range("F33").Value = 80
$ws1.Range("F34").Value = 164
$ws1.Range("F35").Value = 586
$ws1.Range("F36").Value = 457
$ws1.Range("F37").Value = 175
$ws1.Range("F38").Value = 900
$ws1.Range("F39").Value = 690
$ws1.Range("F41").Value = 379
$ws1.Range("F42").Value = 340

# ---- Sheet 2 ----
$ws2.Range("G3").Value = "不可售"
$ws2.Range("F6").Value = 76
$ws2.Range("F12").Value = 50
$ws2.Range("F22").Value = 85
$ws2.Range("F24").Value = 85
$ws2.Range("F25").Value = 427

# ---- Sheet 3 ----
$ws3.Range("F2").Value = 227
$ws3.Range("F3").Value = 2912
$ws3.Range("F6").Value = 290

# ---- Sheet 4 ----
$ws4.Range("F2").Value = 327
$ws4.Range("F3").Value = 227
$ws4.Range("F8").Value = 611
$ws4.Range("C10").Value = "北京·小不点视界-木偶独角戏《千里走单骑》"
$ws4.Range("D10").Value = "南中轴路西侧、永定门以北 天桥艺术中心"
$ws4.Range("E10").Value = "2024.10.02 19:30-10.05 20:30"
$ws4.Range("F10").Value = 0
$ws4.Range("G10").Value = 320
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=92910"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202409/JabgxB8n1727248714888.jpeg"
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = "2024-10-03"
$ws4.Range("C11").Value = "北京·明日方舟同人only-厮守序言"
$ws4.Range("D11").Value = "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
$ws4.Range("E11").Value = "2024.10.03 09:30-10.03 17:00"
$ws4.Range("F11").Value = 453
$ws4.Range("G11").Value = 68
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90959"
$ws4.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202408/rIGY6eyZ1723974119991.jpeg"
$ws4.Range("B12").NumberFormat = "@"
$ws4.Range("B12").Value = "2024-10-04"
$ws4.Range("C12").Value = "北京·第五人格only同人展"
$ws4.Range("D12").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws4.Range("E12").Value = "2024.10.04 10:00-10.04 17:00"
$ws4.Range("F12").Value = 2102
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=89309"
$ws4.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202407/4XsICpa71721046044404.jpeg"
$ws4.Range("C13").Value = "北京·首届SH动漫游戏展"
$ws4.Range("D13").Value = "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
$ws4.Range("E13").Value = "2024.10.04 10:00-10.04 18:00"
$ws4.Range("F13").Value = 845
$ws4.Range("G13").Value = 55
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=91635"
$ws4.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202409/SKe1HMLH1725179765551.jpeg"
$ws4.Range("C14").Value = "帝都·重返未来1999同人ONLY金秋深眠"
$ws4.Range("D14").Value = "华佗路与新源大街交汇处西100米 凯德MALL·大兴"
$ws4.Range("E14").Value = "2024.10.04 10:00-10.05 17:00"
$ws4.Range("F14").Value = 809
$ws4.Range("G14").Value = 68
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=92315"
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202409/YHMYHehz1726129707544.jpeg"
$ws4.Range("B15").NumberFormat = "@"
$ws4.Range("B15").Value = "2024-10-05"
$ws4.Range("C15").Value = "北京·咒术回战同人Only2.0"
$ws4.Range("D15").Value = "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
$ws4.Range("E15").Value = "2024.10.05 09:30-10.05 17:00"
$ws4.Range("F15").Value = 399
$ws4.Range("G15").Value = 65
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=91628"
$ws4.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202408/IsJo7aU61724405528082.jpeg"
$ws4.Range("F16").Value = 63
$ws4.Range("F18").Value = 317
$ws4.Range("F20").Value = 680
$ws4.Range("F21").Value = 25
$ws4.Range("F22").Value = 290
$ws4.Range("F23").Value = 1626
$ws4.Range("F24").Value = 41
$ws4.Range("F25").Value = 50
$ws4.Range("F31").Value = 1450
$ws4.Range("F34").Value = 515
$ws4.Range("F35").Value = 566
$ws4.Range("F36").Value = 406
$ws4.Range("F38").Value = 2269
$ws4.Range("F39").Value = 80
$ws4.Range("F40").Value = 164
$ws4.Range("F41").Value = 586
$ws4.Range("F42").Value = 457
$ws4.Range("F43").Value = 175
$ws4.Range("F44").Value = 900
$ws4.Range("F46").Value = 85
$ws4.Range("F47").Value = 427
$ws4.Range("F48").Value = 690
